$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.799319333333333
$ws.Range("H2").Value = 11.397958
$ws.Range("I2").Value = 0.04516380397110807
$ws.Range("J2").Value = 0.04516380397110807
$ws.Range("M2").Value = 51.15371566666666
$ws.Range("N2").Value = 153.461147
$ws.Range("O2").Value = 0.3311207986511828
$ws.Range("P2").Value = 0.3311207986511828
$ws.Range("Q2").Value = 194.3493009042029
$ws.Range("R2").Value = 1749.143708137826
$ws.Range("S2").Value = 0.01495467484103877
$ws.Range("T2").Value = 0.01495467484103877
$ws.Range("G3").Value = 3.799319333333333
$ws.Range("H3").Value = 11.397958
$ws.Range("I3").Value = 0.04516380397110807
$ws.Range("J3").Value = 0.04516380397110807
$ws.Range("M3").Value = 53.36146666666667
$ws.Range("O3").Value = 0.3454116915964105
$ws.Range("P3").Value = 0.3454116915964106
$ws.Range("Q3").Value = 202.7372519616889
$ws.Range("R3").Value = 1824.6352676552
$ws.Range("S3").Value = 0.01560010592858912
$ws.Range("T3").Value = 0.01560010592858912
$ws.Range("G4").Value = 3.799319333333333
$ws.Range("H4").Value = 11.397958
$ws.Range("I4").Value = 0.04516380397110807
$ws.Range("J4").Value = 0.04516380397110807
$ws.Range("M4").Value = 25.00653133333333
$ws.Range("N4").Value = 75.019594
$ws.Range("O4").Value = 0.1618686447050176
$ws.Range("P4").Value = 0.1618686447050176
$ws.Range("Q4").Value = 95.00779795433911
$ws.Range("R4").Value = 855.070181589052
$ws.Range("S4").Value = 0.007310603738526357
$ws.Range("T4").Value = 0.007310603738526357
$ws.Range("G5").Value = 3.799319333333333
$ws.Range("H5").Value = 11.397958
$ws.Range("I5").Value = 0.04516380397110807
$ws.Range("J5").Value = 0.04516380397110807
$ws.Range("M5").Value = 6.481347
$ws.Range("N5").Value = 19.444041
$ws.Range("O5").Value = 0.0419541135381084
$ws.Range("P5").Value = 0.0419541135381084
$ws.Range("Q5").Value = 24.624706963142
$ws.Range("R5").Value = 221.622362668278
$ws.Range("S5").Value = 0.001894807359616739
$ws.Range("T5").Value = 0.001894807359616739
$ws.Range("G6").Value = 3.799319333333333
$ws.Range("H6").Value = 11.397958
$ws.Range("I6").Value = 0.04516380397110807
$ws.Range("J6").Value = 0.04516380397110807
$ws.Range("M6").Value = 18.483507
$ws.Range("N6").Value = 55.450521
$ws.Range("O6").Value = 0.1196447515092806
$ws.Range("P6").Value = 0.1196447515092806
$ws.Range("Q6").Value = 70.224745492902
$ws.Range("R6").Value = 632.022709436118
$ws.Range("S6").Value = 0.005403612103337086
$ws.Range("T6").Value = 0.005403612103337087
$ws.Range("I7").Value = 0.7704071060399821
$ws.Range("J7").Value = 0.770407106039982
$ws.Range("M7").Value = 51.15371566666666
$ws.Range("N7").Value = 153.461147
$ws.Range("O7").Value = 0.3311207986511828
$ws.Range("P7").Value = 0.3311207986511828
$ws.Range("Q7").Value = 3315.223017226002
$ws.Range("R7").Value = 29837.00715503402
$ws.Range("S7").Value = 0.2550978162385054
$ws.Range("T7").Value = 0.2550978162385053
$ws.Range("I8").Value = 0.7704071060399821
$ws.Range("J8").Value = 0.770407106039982
$ws.Range("M8").Value = 53.36146666666667
$ws.Range("O8").Value = 0.3454116915964105
$ws.Range("P8").Value = 0.3454116915964106
$ws.Range("Q8").Value = 3458.305231999956
$ws.Range("R8").Value = 31124.7470879996
$ws.Range("S8").Value = 0.2661076217151654
$ws.Range("T8").Value = 0.2661076217151654
$ws.Range("I9").Value = 0.7704071060399821
$ws.Range("J9").Value = 0.770407106039982
$ws.Range("M9").Value = 25.00653133333333
$ws.Range("N9").Value = 75.019594
$ws.Range("O9").Value = 0.1618686447050176
$ws.Range("P9").Value = 0.1618686447050176
$ws.Range("Q9").Value = 1620.649197752638
$ws.Range("R9").Value = 14585.84277977375
$ws.Range("S9").Value = 0.1247047541258067
$ws.Range("T9").Value = 0.1247047541258067
$ws.Range("I10").Value = 0.7704071060399821
$ws.Range("J10").Value = 0.770407106039982
$ws.Range("M10").Value = 6.481347
$ws.Range("N10").Value = 19.444041
$ws.Range("O10").Value = 0.0419541135381084
$ws.Range("P10").Value = 0.0419541135381084
$ws.Range("Q10").Value = 420.049853211941
$ws.Range("R10").Value = 3780.448678907469
$ws.Range("S10").Value = 0.03232174719736693
$ws.Range("T10").Value = 0.03232174719736692
$ws.Range("I11").Value = 0.7704071060399821
$ws.Range("J11").Value = 0.770407106039982
$ws.Range("M11").Value = 18.483507
$ws.Range("N11").Value = 55.450521
$ws.Range("O11").Value = 0.1196447515092806
$ws.Range("P11").Value = 0.1196447515092806
$ws.Range("Q11").Value = 1197.898276730421
$ws.Range("R11").Value = 10781.08449057379
$ws.Range("S11").Value = 0.09217516676313765
$ws.Range("T11").Value = 0.09217516676313765
$ws.Range("G12").Value = 2.946166666666667
$ws.Range("H12").Value = 8.8385
$ws.Range("I12").Value = 0.03502208741238024
$ws.Range("J12").Value = 0.03502208741238024
$ws.Range("M12").Value = 51.15371566666666
$ws.Range("N12").Value = 153.461147
$ws.Range("O12").Value = 0.3311207986511828
$ws.Range("P12").Value = 0.3311207986511828
$ws.Range("Q12").Value = 150.7073719732778
$ws.Range("R12").Value = 1356.3663477595
$ws.Range("S12").Value = 0.01159654155441888
$ws.Range("T12").Value = 0.01159654155441888
$ws.Range("G13").Value = 2.946166666666667
$ws.Range("H13").Value = 8.8385
$ws.Range("I13").Value = 0.03502208741238024
$ws.Range("J13").Value = 0.03502208741238024
$ws.Range("M13").Value = 53.36146666666667
$ws.Range("O13").Value = 0.3454116915964105
$ws.Range("P13").Value = 0.3454116915964106
$ws.Range("Q13").Value = 157.2117743777778
$ws.Range("R13").Value = 1414.9059694
$ws.Range("S13").Value = 0.01209703845634761
$ws.Range("T13").Value = 0.01209703845634762
$ws.Range("G14").Value = 2.946166666666667
$ws.Range("H14").Value = 8.8385
$ws.Range("I14").Value = 0.03502208741238024
$ws.Range("J14").Value = 0.03502208741238024
$ws.Range("M14").Value = 25.00653133333333
$ws.Range("N14").Value = 75.019594
$ws.Range("O14").Value = 0.1618686447050176
$ws.Range("P14").Value = 0.1618686447050176
$ws.Range("Q14").Value = 73.67340906322221
$ws.Range("R14").Value = 663.0606815689999
$ws.Range("S14").Value = 0.005668977824182648
$ws.Range("T14").Value = 0.005668977824182648
$ws.Range("G15").Value = 2.946166666666667
$ws.Range("H15").Value = 8.8385
$ws.Range("I15").Value = 0.03502208741238024
$ws.Range("J15").Value = 0.03502208741238024
$ws.Range("M15").Value = 6.481347
$ws.Range("N15").Value = 19.444041
$ws.Range("O15").Value = 0.0419541135381084
$ws.Range("P15").Value = 0.0419541135381084
$ws.Range("Q15").Value = 19.0951284865
$ws.Range("R15").Value = 171.8561563785
$ws.Range("S15").Value = 0.001469320631640558
$ws.Range("T15").Value = 0.001469320631640558
$ws.Range("G16").Value = 2.946166666666667
$ws.Range("H16").Value = 8.8385
$ws.Range("I16").Value = 0.03502208741238024
$ws.Range("J16").Value = 0.03502208741238024
$ws.Range("M16").Value = 18.483507
$ws.Range("N16").Value = 55.450521
$ws.Range("O16").Value = 0.1196447515092806
$ws.Range("P16").Value = 0.1196447515092806
$ws.Range("Q16").Value = 54.4554922065
$ws.Range("R16").Value = 490.0994298585
$ws.Range("S16").Value = 0.004190208945790539
$ws.Range("T16").Value = 0.004190208945790539
$ws.Range("G17").Value = 10.035916
$ws.Range("H17").Value = 30.107748
$ws.Range("I17").Value = 0.1193003543865946
$ws.Range("J17").Value = 0.1193003543865946
$ws.Range("M17").Value = 51.15371566666666
$ws.Range("N17").Value = 153.461147
$ws.Range("O17").Value = 0.3311207986511828
$ws.Range("P17").Value = 0.3311207986511828
$ws.Range("Q17").Value = 513.3743935185506
$ws.Range("R17").Value = 4620.369541666955
$ws.Range("S17").Value = 0.03950282862385835
$ws.Range("T17").Value = 0.03950282862385834
$ws.Range("G18").Value = 10.035916
$ws.Range("H18").Value = 30.107748
$ws.Range("I18").Value = 0.1193003543865946
$ws.Range("J18").Value = 0.1193003543865946
$ws.Range("M18").Value = 53.36146666666667
$ws.Range("O18").Value = 0.3454116915964105
$ws.Range("P18").Value = 0.3454116915964106
$ws.Range("Q18").Value = 535.5311971034668
$ws.Range("R18").Value = 4819.7807739312
$ws.Range("S18").Value = 0.0412077372167249
$ws.Range("T18").Value = 0.0412077372167249
$ws.Range("G19").Value = 10.035916
$ws.Range("H19").Value = 30.107748
$ws.Range("I19").Value = 0.1193003543865946
$ws.Range("J19").Value = 0.1193003543865946
$ws.Range("M19").Value = 25.00653133333333
$ws.Range("N19").Value = 75.019594
$ws.Range("O19").Value = 0.1618686447050176
$ws.Range("P19").Value = 0.1618686447050176
$ws.Range("Q19").Value = 250.9634479127013
$ws.Range("R19").Value = 2258.671031214312
$ws.Range("S19").Value = 0.01931098667738637
$ws.Range("T19").Value = 0.01931098667738637
$ws.Range("G20").Value = 10.035916
$ws.Range("H20").Value = 30.107748
$ws.Range("I20").Value = 0.1193003543865946
$ws.Range("J20").Value = 0.1193003543865946
$ws.Range("M20").Value = 6.481347
$ws.Range("N20").Value = 19.444041
$ws.Range("O20").Value = 0.0419541135381084
$ws.Range("P20").Value = 0.0419541135381084
$ws.Range("Q20").Value = 65.046254058852
$ws.Range("R20").Value = 585.416286529668
$ws.Range("S20").Value = 0.005005140613071759
$ws.Range("T20").Value = 0.005005140613071758
$ws.Range("G21").Value = 10.035916
$ws.Range("H21").Value = 30.107748
$ws.Range("I21").Value = 0.1193003543865946
$ws.Range("J21").Value = 0.1193003543865946
$ws.Range("M21").Value = 18.483507
$ws.Range("N21").Value = 55.450521
$ws.Range("O21").Value = 0.1196447515092806
$ws.Range("P21").Value = 0.1196447515092806
$ws.Range("Q21").Value = 185.498923637412
$ws.Range("R21").Value = 1669.490312736708
$ws.Range("S21").Value = 0.01427366125555323
$ws.Range("T21").Value = 0.01427366125555323
$ws.Range("G22").Value = 2.532664666666667
$ws.Range("H22").Value = 7.597994
$ws.Range("I22").Value = 0.03010664818993502
$ws.Range("J22").Value = 0.03010664818993501
$ws.Range("M22").Value = 51.15371566666666
$ws.Range("N22").Value = 153.461147
$ws.Range("O22").Value = 0.3311207986511828
$ws.Range("P22").Value = 0.3311207986511828
$ws.Range("Q22").Value = 129.5552082376798
$ws.Range("R22").Value = 1165.996874139118
$ws.Range("S22").Value = 0.009968937393361471
$ws.Range("T22").Value = 0.009968937393361469
$ws.Range("G23").Value = 2.532664666666667
$ws.Range("H23").Value = 7.597994
$ws.Range("I23").Value = 0.03010664818993502
$ws.Range("J23").Value = 0.03010664818993501
$ws.Range("M23").Value = 53.36146666666667
$ws.Range("O23").Value = 0.3454116915964105
$ws.Range("P23").Value = 0.3454116915964106
$ws.Range("Q23").Value = 135.1467011881778
$ws.Range("R23").Value = 1216.3203106936
$ws.Range("S23").Value = 0.01039918827958346
$ws.Range("T23").Value = 0.01039918827958346
$ws.Range("G24").Value = 2.532664666666667
$ws.Range("H24").Value = 7.597994
$ws.Range("I24").Value = 0.03010664818993502
$ws.Range("J24").Value = 0.03010664818993501
$ws.Range("M24").Value = 25.00653133333333
$ws.Range("N24").Value = 75.019594
$ws.Range("O24").Value = 0.1618686447050176
$ws.Range("P24").Value = 0.1618686447050176
$ws.Range("Q24").Value = 63.33315834382622
$ws.Range("R24").Value = 569.998425094436
$ws.Range("S24").Value = 0.004873322339115554
$ws.Range("T24").Value = 0.004873322339115553
$ws.Range("G25").Value = 2.532664666666667
$ws.Range("H25").Value = 7.597994
$ws.Range("I25").Value = 0.03010664818993502
$ws.Range("J25").Value = 0.03010664818993501
$ws.Range("M25").Value = 6.481347
$ws.Range("N25").Value = 19.444041
$ws.Range("O25").Value = 0.0419541135381084
$ws.Range("P25").Value = 0.0419541135381084
$ws.Range("Q25").Value = 16.415078539306
$ws.Range("R25").Value = 147.735706853754
$ws.Range("S25").Value = 0.001263097736412419
$ws.Range("T25").Value = 0.001263097736412419
$ws.Range("G26").Value = 2.532664666666667
$ws.Range("H26").Value = 7.597994
$ws.Range("I26").Value = 0.03010664818993502
$ws.Range("J26").Value = 0.03010664818993501
$ws.Range("M26").Value = 18.483507
$ws.Range("N26").Value = 55.450521
$ws.Range("O26").Value = 0.1196447515092806
$ws.Range("P26").Value = 0.1196447515092806
$ws.Range("Q26").Value = 70.224745492902
$ws.Range("R26").Value = 421.312725854874
$ws.Range("S26").Value = 0.003602102441462108
$ws.Range("T26").Value = 0.003602102441462108
